# New crime data collected - weekly CompStat report refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: volume/report-week labels (rich-text shared strings flatten to
# plain text through this COM layer, but the run font == the cell's own
# style font, so the rendered result is identical either way).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/6/2025  Through  1/12/2025"

# ---------------------------------------------------------------------------
# Cells that flip from the "0" / "***.*" text placeholder to a real number.
# Copy format+value from a stable numeric template cell, then overwrite the
# value so the resulting style id matches what the workbook already uses for
# plain numbers (14) / one-decimal numbers (15).
# ---------------------------------------------------------------------------
$numInt = $ws.Range("G16")    # template cell: style 14 (integer format)
$numDec = $ws.Range("L16")    # template cell: style 15 (one-decimal format)

function Set-NumFromText($addr, $template, $value) {
    $template.Copy($ws.Range($addr))
    $ws.Range($addr).Value = $value
}

Set-NumFromText "M16" $numDec 0
Set-NumFromText "D17" $numInt 1
Set-NumFromText "E17" $numDec 200
Set-NumFromText "J17" $numInt 1
Set-NumFromText "K17" $numDec 400
Set-NumFromText "D18" $numInt 1
Set-NumFromText "E18" $numDec -100
Set-NumFromText "J18" $numInt 1
Set-NumFromText "K18" $numDec -100
Set-NumFromText "L18" $numDec -100
Set-NumFromText "M18" $numDec -100
Set-NumFromText "C19" $numInt 1
Set-NumFromText "I19" $numInt 1
Set-NumFromText "C20" $numInt 1
Set-NumFromText "F20" $numInt 1
Set-NumFromText "I20" $numInt 1
Set-NumFromText "I26" $numInt 3
Set-NumFromText "L28" $numDec -100

# ---------------------------------------------------------------------------
# Cells that flip from a real number to the "0" / "***.*" text placeholder.
# Copy format+value straight from a stable cell that already holds the
# placeholder text so the existing shared-string entries (20 / 21) are reused
# verbatim instead of minting new strings.
# ---------------------------------------------------------------------------
$text0 = $ws.Range("D16")     # template cell: style 13, shared string 20 ("0")
$textStar = $ws.Range("E16")  # template cell: style 13, shared string 21 ("***.*")

$text0.Copy($ws.Range("C16"))
$text0.Copy($ws.Range("C23"))
$text0.Copy($ws.Range("D28"))
$text0.Copy($ws.Range("G33"))

$textStar.Copy($ws.Range("E28"))
$textStar.Copy($ws.Range("H33"))

# ---------------------------------------------------------------------------
# Plain numeric value updates (style unchanged).
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = 5
$ws.Range("H16").Value = 400
$ws.Range("N16").Value = -85.714285714285

$ws.Range("C17").Value = 3
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 11.111111111111
$ws.Range("I17").Value = 5
$ws.Range("M17").Value = 66.666666666666
$ws.Range("N17").Value = 25

$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = -66.666666666666

$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -75
$ws.Range("F19").Value = 4
$ws.Range("G19").Value = 6
$ws.Range("H19").Value = -33.333333333333
$ws.Range("J19").Value = 5
$ws.Range("K19").Value = -80
$ws.Range("L19").Value = -85.714285714285
$ws.Range("M19").Value = -83.333333333333
$ws.Range("N19").Value = -85.714285714285

$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("H20").Value = -90
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = -50
$ws.Range("N20").Value = -90.909090909090

$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 7
$ws.Range("E21").Value = -28.571428571428
$ws.Range("F21").Value = 21
$ws.Range("G21").Value = 29
$ws.Range("H21").Value = -27.586206896551
$ws.Range("I21").Value = 8
$ws.Range("J21").Value = 9
$ws.Range("K21").Value = -11.111111111111
$ws.Range("L21").Value = -46.666666666666
$ws.Range("M21").Value = -27.272727272727
$ws.Range("N21").Value = -80.952380952380

$ws.Range("F23").Value = 3
$ws.Range("I23").Value = 2
$ws.Range("M23").Value = 0

$ws.Range("C24").Value = 7
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 28
$ws.Range("G24").Value = 38
$ws.Range("H24").Value = -26.315789473684
$ws.Range("I24").Value = 10
$ws.Range("J24").Value = 13
$ws.Range("K24").Value = -23.076923076923
$ws.Range("L24").Value = -61.538461538461
$ws.Range("M24").Value = -9.090909090909

$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 6
$ws.Range("H25").Value = -57.142857142857
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 5
$ws.Range("K25").Value = -60
$ws.Range("L25").Value = -75

$ws.Range("C26").Value = 3
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 12
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = -7.692307692307
$ws.Range("J26").Value = 6
$ws.Range("K26").Value = -50
$ws.Range("L26").Value = -25
$ws.Range("M26").Value = -62.5

$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -50

$ws.Range("J39").Value = 2
$ws.Range("K39").Value = 100
$ws.Range("M39").Value = -60
$ws.Range("N39").Value = -71.428571428571

$ws.Range("J46").Value = 481
$ws.Range("K46").Value = 7.606263982102
$ws.Range("L46").Value = -27.777777777777
$ws.Range("M46").Value = -66.689750692520
$ws.Range("N46").Value = -69.576217583807
